# Applies the weekly fruit/vegetable data refresh: the price rows for
# "Fruta, Vega Modelo de Temuco - Frambuesa" were rotated so that each
# row's Fecha/Volumen/Precio/Origen block now corresponds to a different
# week than before. Rows 6 and 8 are untouched; the remaining rows
# (2,3,4,5,7,9,10,11) receive the values previously held by another row,
# per the mapping: new row <- old row
#   2<-3, 3<-5, 4<-7, 5<-4, 7<-2, 9<-11, 10<-9, 11<-10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) values for the columns that move, keyed
# by row number, before any writes happen.
$cols = @("D", "M", "N", "O", "P", "R", "S")
$orig = @{}
foreach ($r in 2..11) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $orig[$r] = $rowVals
}

# Mapping of destination row -> source row (source row's old values are
# written into the destination row).
$map = @{
    2  = 3
    3  = 5
    4  = 7
    5  = 4
    7  = 2
    9  = 11
    10 = 9
    11 = 10
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $orig[$srcRow][$c]
    }
}
